# Applies a row-content permutation within several blocks of the "Artfynd"
# sheet. The underlying observation records were re-ordered (rows keep their
# row number / position, but the data that lives in each row moves to a
# different row), while row 1 (headers) and all unrelated rows are untouched.
#
# For each block below, the map is: target row -> source row (i.e. the row
# that currently -- before this script runs -- holds the data that should
# end up in the target row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that hold plain numbers (copied with Value2 / Value, no special
# formatting needed).
$numericCols = @(1, 2, 5, 17, 18, 19)          # A, B, E, Q, R, S

# Columns that hold TRUE/FALSE booleans.
$boolCols = @(30, 31, 33)                       # AD, AE, AG

# Columns that hold text -- including values that look like dates/times
# ("2023-08-31", "00:00") which Excel would otherwise auto-convert, so the
# destination cell's format is forced to Text ("@") before the value is
# written. (Columns I, AT, AY are intentionally excluded: they are blank on
# every row in every affected block, both before and after, so there is
# nothing to move for them.)
$textCols = @(3, 4, 6, 7, 8, 16, 20, 21, 22, 23, 25, 26, 27, 28, 49, 50)
# C, D, F, G, H, P, T, U, V, W, Y, Z, AA, AB, AW, AX

# The "Bestamningsmetod" column -- present (but always empty) on only some
# rows of each block; its presence must travel with the rest of the row's
# data. The cell is always blank, so there's nothing distinguishable to read
# back at runtime (an empty cell and a missing one read identically through
# the object model) -- instead we record, from the original file, which row
# numbers currently carry that empty AF cell.
$afCol = 32                                     # AF
$afPresentRows = @(14, 24, 26, 34, 39, 52, 53, 80, 89)

function Read-RowData($row) {
    $data = @{}
    foreach ($c in $numericCols) { $data[$c] = $ws.Cells.Item($row, $c).Value2 }
    foreach ($c in $boolCols)    { $data[$c] = $ws.Cells.Item($row, $c).Value2 }
    foreach ($c in $textCols)    { $data[$c] = $ws.Cells.Item($row, $c).Value2 }
    $afCell = $ws.Cells.Item($row, $afCol)
    $data[$afCol] = $afCell.Value2
    $data["__af_present__"] = ($afPresentRows -contains $row)
    return $data
}

function Write-RowData($row, $data) {
    foreach ($c in $numericCols) {
        $v = $data[$c]
        $cell = $ws.Cells.Item($row, $c)
        if ($v -eq $null) { $cell.ClearContents() } else { $cell.Value = $v }
    }
    foreach ($c in $boolCols) {
        $cell = $ws.Cells.Item($row, $c)
        $cell.Value = $data[$c]
    }
    foreach ($c in $textCols) {
        $cell = $ws.Cells.Item($row, $c)
        $v = $data[$c]
        # Force Text format so date-/time-looking strings ("2023-08-31",
        # "00:00") are not auto-converted to date/time serials on write.
        $cell.NumberFormat = "@"
        if ($v -ne $null -and $v -ne "") {
            $cell.Value = $v
        } elseif ($cell.Value2 -ne $null -and $cell.Value2 -ne "") {
            # Source is blank but destination currently isn't -- clear it.
            $cell.ClearContents()
            $cell.NumberFormat = "@"
        }
    }
    $afCell = $ws.Cells.Item($row, $afCol)
    if ($data["__af_present__"]) {
        # AF is blank on every row that carries it -- setting the format
        # (without touching .Value) is enough to make the engine persist it
        # as a present-but-empty cell, mirroring the source file.
        $afCell.NumberFormat = "@"
    } else {
        $afCell.ClearContents()
    }
}

function Apply-Block($map) {
    # Read every source row's data BEFORE writing anything, since several
    # source rows are themselves targets within the same block.
    $snapshots = @{}
    foreach ($targetRow in $map.Keys) {
        $sourceRow = $map[$targetRow]
        $snapshots[$targetRow] = Read-RowData $sourceRow
    }
    foreach ($targetRow in $map.Keys) {
        Write-RowData $targetRow $snapshots[$targetRow]
    }
}

# Block 1: rows 12-16
$block1 = @{
    12 = 13
    13 = 16
    14 = 12
    15 = 14
    16 = 15
}
Apply-Block $block1

# Block 2: rows 26-28
$block2 = @{
    26 = 28
    27 = 26
    28 = 27
}
Apply-Block $block2

# Block 3: rows 51-54
$block3 = @{
    51 = 52
    52 = 53
    53 = 54
    54 = 51
}
Apply-Block $block3

# Block 4: rows 76-80
$block4 = @{
    76 = 78
    77 = 80
    78 = 76
    79 = 77
    80 = 79
}
Apply-Block $block4
